$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 7550.3335
$ws.Range("I80").Value = 180
$ws.Range("J80").Value = 10230.454
$ws.Range("K80").Value = 540
$ws.Range("L80").Value = 30691.362
$ws.Range("M80").Value = 458
$ws.Range("N80").Value = -32687.362
$ws.Range("H83").Value = 7550.3335
$ws.Range("I83").Value = 180
$ws.Range("J83").Value = 10230.454
$ws.Range("K83").Value = 1620
$ws.Range("L83").Value = 92074.086
$ws.Range("M83").Value = 3372
$ws.Range("N83").Value = -102058.086
$ws.Range("H115").Value = 4498
$ws.Range("I115").Value = 1280
$ws.Range("J115").Value = 6643.3335
$ws.Range("K115").Value = 3840
$ws.Range("L115").Value = 19930.0005
$ws.Range("M115").Value = -2273
$ws.Range("N115").Value = -23064.0005
$ws.Range("H125").Value = 3451.2307
$ws.Range("I125").Value = 1915.5
$ws.Range("J125").Value = 3730.4546
$ws.Range("K125").Value = 17239.5
$ws.Range("L125").Value = 33574.0914
$ws.Range("M125").Value = -14779.5
$ws.Range("N125").Value = -38494.0914
$ws.Range("H127").Value = 2407.7
$ws.Range("I127").Value = 838.5
$ws.Range("K127").Value = 2515.5
$ws.Range("M127").Value = 2444.5
$ws.Range("H129").Value = 1127.4507
$ws.Range("J129").Value = 1137.9857
$ws.Range("L129").Value = 3413.9571
$ws.Range("N129").Value = -13413.9571
$ws.Range("H138").Value = 3554.2646
$ws.Range("I138").Value = 2743.0908
$ws.Range("J138").Value = 3710.8071
$ws.Range("K138").Value = 8229.2724
$ws.Range("L138").Value = 11132.4213
$ws.Range("M138").Value = -3089.2724
$ws.Range("N138").Value = -21412.4213
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 50000
$ws.Range("K3").Value = 50000
$ws.Range("M3").Value = -49885
$ws.Range("H22").Value = 889.5294
$ws.Range("I22").Value = 740.1539
$ws.Range("J22").Value = 1375
$ws.Range("K22").Value = 740.1539
$ws.Range("L22").Value = 1375
$ws.Range("M22").Value = -441.1539
$ws.Range("N22").Value = -1973
$ws.Range("H26").Value = 54000
$ws.Range("I26").Value = 54000
$ws.Range("K26").Value = 54000
$ws.Range("M26").Value = -53670
$ws.Range("H32").Value = 16621.695
$ws.Range("I32").Value = 20532.266
$ws.Range("J32").Value = 5541.75
$ws.Range("K32").Value = 20532.266
$ws.Range("L32").Value = 5541.75
$ws.Range("M32").Value = -20245.266
$ws.Range("N32").Value = -6115.75
$ws.Range("H39").Value = 16200
$ws.Range("I39").Value = 16200
$ws.Range("K39").Value = 16200
$ws.Range("M39").Value = -15680
$ws.Range("H41").Value = 2133.8667
$ws.Range("I41").Value = 1929.1428
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 1929.1428
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -1515.1428
$ws.Range("N41").Value = -5828
$ws.Range("H122").Value = 7149.3335
$ws.Range("I122").Value = 7452.2354
$ws.Range("K122").Value = 22356.7062
$ws.Range("M122").Value = -19906.7062
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 9866.666999999999
$ws.Range("I8").Value = 6236.364
$ws.Range("K8").Value = 6236.364
$ws.Range("M8").Value = -6096.364
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2060.9275
$ws.Range("I31").Value = 1928.2609
$ws.Range("J31").Value = 2326.261
$ws.Range("K31").Value = 1928.2609
$ws.Range("L31").Value = 2326.261
$ws.Range("M31").Value = -1633.2609
$ws.Range("N31").Value = -2916.261
$ws.Range("H34").Value = 2060.9275
$ws.Range("I34").Value = 1928.2609
$ws.Range("J34").Value = 2326.261
$ws.Range("K34").Value = 1928.2609
$ws.Range("L34").Value = 2326.261
$ws.Range("M34").Value = -1726.2609
$ws.Range("N34").Value = -2730.261
$ws.Range("H42").Value = 18666.666
$ws.Range("I42").Value = 25500
$ws.Range("J42").Value = 5000
$ws.Range("K42").Value = 25500
$ws.Range("L42").Value = 5000
$ws.Range("M42").Value = -24907
$ws.Range("N42").Value = -6186
$ws.Range("H62").Value = 128201.25
$ws.Range("I62").Value = 128201.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 128201.25
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -127577.25
$ws.Range("H65").Value = 128201.25
$ws.Range("I65").Value = 128201.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 641006.25
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -637886.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 178533.33
$ws.Range("I4").Value = 250300
$ws.Range("J4").Value = 35000
$ws.Range("K4").Value = 750900
$ws.Range("L4").Value = 105000
$ws.Range("M4").Value = -750788
$ws.Range("N4").Value = -105224
$ws.Range("H29").Value = 1096.6666
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1096.6666
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 3289.9998
$ws.Range("N29").Value = -3843.9998
$ws.Range("H68").Value = 143671.11
$ws.Range("I68").Value = 257002.52
$ws.Range("J68").Value = 1092.9032
$ws.Range("K68").Value = 771007.5599999999
$ws.Range("L68").Value = 3278.7096
$ws.Range("M68").Value = -770196.5599999999
$ws.Range("N68").Value = -4900.7096
$ws.Range("H71").Value = 143671.11
$ws.Range("I71").Value = 257002.52
$ws.Range("J71").Value = 1092.9032
$ws.Range("K71").Value = 2313022.68
$ws.Range("L71").Value = 9836.1288
$ws.Range("M71").Value = -2308966.68
$ws.Range("N71").Value = -17948.1288
$ws.Range("H103").Value = 2504
$ws.Range("J103").Value = 2465.6
$ws.Range("L103").Value = 7396.799999999999
$ws.Range("N103").Value = -9154.799999999999
$ws.Range("H113").Value = 213415.48
$ws.Range("I113").Value = 370938.88
$ws.Range("J113").Value = 758.9
$ws.Range("K113").Value = 1112816.64
$ws.Range("L113").Value = 2276.7
$ws.Range("M113").Value = -1110646.64
$ws.Range("N113").Value = -6616.7
$ws.Range("H131").Value = 20848.568
$ws.Range("I131").Value = 1310.4166
$ws.Range("J131").Value = 26860.309
$ws.Range("K131").Value = 3931.2498
$ws.Range("L131").Value = 80580.927
$ws.Range("M131").Value = 1108.7502
$ws.Range("N131").Value = -90660.927
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28999.5
$ws.Range("J39").Value = 28999.5
$ws.Range("L39").Value = 28999.5
$ws.Range("N39").Value = -30063.5
$ws.Range("H70").Value = 8410.736999999999
$ws.Range("I70").Value = 8775.25
$ws.Range("J70").Value = 6466.6665
$ws.Range("K70").Value = 8775.25
$ws.Range("L70").Value = 6466.6665
$ws.Range("M70").Value = -8505.25
$ws.Range("N70").Value = -7006.6665
$ws.Range("H73").Value = 8410.736999999999
$ws.Range("I73").Value = 8775.25
$ws.Range("J73").Value = 6466.6665
$ws.Range("K73").Value = 8775.25
$ws.Range("L73").Value = 6466.6665
$ws.Range("M73").Value = -7839.25
$ws.Range("N73").Value = -8338.666499999999
$ws.Range("H123").Value = 8823.739
$ws.Range("J123").Value = 8823.739
$ws.Range("L123").Value = 8823.739
$ws.Range("N123").Value = -13723.739
$ws.Range("H132").Value = 3393.9048
$ws.Range("I132").Value = 2384.1667
$ws.Range("J132").Value = 4740.222
$ws.Range("K132").Value = 7152.500100000001
$ws.Range("L132").Value = 14220.666
$ws.Range("M132").Value = -4622.500100000001
$ws.Range("N132").Value = -19280.666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 31059.334
$ws.Range("J98").Value = 31059.334
$ws.Range("L98").Value = 31059.334
$ws.Range("N98").Value = -37049.334
$ws.Range("H119").Value = 34333.332
$ws.Range("J119").Value = 34333.332
$ws.Range("L119").Value = 34333.332
$ws.Range("N119").Value = -44009.332
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 50000
$ws.Range("I22").Value = 50000
$ws.Range("K22").Value = 50000
$ws.Range("M22").Value = -49707
